$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header info
$ws.Range("C2").Value = "Hartmut"

# Card number must stay as text (avoid numeric precision loss / auto
# number conversion on a 16-digit value) - use the classic text-prefix,
# then re-apply the original (non quote-prefixed) number format from the
# cell above so the cell style index is unaffected by the prefix trick.
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Value = "Mohaupt"

# Opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 27.02.2025"

# Transaction row 6
$ws.Range("B6").Value = "02.03."
$ws.Range("C6").Value = "03.03."
$ws.Range("D6").Value = "BEITRAG Allianz SE K-4097343"
$ws.Range("E6").Value = "55,98-"

# Transaction row 7
$ws.Range("B7").Value = "05.03."
$ws.Range("C7").Value = "06.03."
$ws.Range("D7").Value = "RECHNUNG VODAFONE GMBH 86434328"
$ws.Range("E7").Value = "38,10-"

# Transaction row 8
$ws.Range("B8").Value = "06.03."
$ws.Range("C8").Value = "07.03."
$ws.Range("D8").Value = "KARTENZ./06.03 REWE RO"
$ws.Range("E8").Value = "63,80-"

# Transaction row 9 - previously blank, now a new transaction. Copy the
# formatting from row 8's cells first so the new row matches the rest of
# the statement's look (right-aligned amount etc.), then fill in values.
$ws.Range("B8:E8").Copy() | Out-Null
$ws.Range("B9:E9").PasteSpecial(-4122) | Out-Null
$ws.Range("B9").Value = "07.03."
$ws.Range("C9").Value = "08.03."
$ws.Range("D9").Value = "AMAZON.DE MKTPLC EU SZEJOQ"
$ws.Range("E9").Value = "148,11-"

# Closing balance line
$ws.Range("D12").Value = "KONTOSTAND AM 09.03.2025"
$ws.Range("E12").Value = "305,99-"

# Next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 16.03.2025"
